$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of the k column (J)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Row 14-17: summary statistics with labels in column A and values in column B
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"

$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Bold, size-12 font with vertically-centered alignment for the summary values.
# Apply directly to B14, then format-paint the same style onto B15:B17 so the
# workbook doesn't accumulate transient intermediate cell-format records.
$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108
$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)

# Selection matching the target state
$ws.Range("A14:B17").Select() | Out-Null

# Page setup (paper size / orientation) as captured by the resave
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
